$wb = $excel.ActiveWorkbook

# Add the new "Fouling Data" worksheet after the last existing sheet ("data")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Fouling Data"

# Populate headers and data
$newSheet.Range("A1").Value = "WaterFouling"
$newSheet.Range("B1").Value = "ChemicalFouling"
$newSheet.Range("A2").Value = "Rhine"
$newSheet.Range("B2").Value = "halogenated alkanes"

# Header style - bold + centered
$newSheet.Range("A1:B1").Font.Bold = $true
$newSheet.Range("A1:B1").HorizontalAlignment = -4108

# Make this sheet the active tab
$newSheet.Activate()
